$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 23:52"

# --- Re-sorted pairs: swap the country label between two adjacent rows ---
# (the underlying counts were refreshed and the two countries swapped rank order)
$ws.Range("A31").Value = "Pakistan"
$ws.Range("A32").Value = "Mexico"

$ws.Range("A50").Value = "Colombia"
$ws.Range("A51").Value = "Banglades"

$ws.Range("A132").Value = "Gabon"
$ws.Range("A133").Value = "Martinica"

# --- Updated numeric data (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$rowValues = @(916348, 29906, 93275, 771331, 14932, 1508, 51742)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

# Row 31: Pakistan (new rank 1 of the swapped pair)
$rowValues = @(11940, 883, 2755, 8932, 111, 18, 253)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(31, $col).Value = $v
    $col = $col + 1
}

# Row 32: Mexico (new rank 2 of the swapped pair)
$rowValues = @(11633, 1089, 2627, 7937, 378, 99, 1069)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(32, $col).Value = $v
    $col = $col + 1
}

# Row 50: Colombia
$rowValues = @(4881, 320, 1003, 3653, 98, 10, 225)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(50, $col).Value = $v
    $col = $col + 1
}

# Row 51: Banglades
$rowValues = @(4689, 503, 112, 4446, 1, 4, 131)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(51, $col).Value = $v
    $col = $col + 1
}

# Row 62: Barein
$rowValues = @(2518, 301, 1113, 1397, 1, 0, 8)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(62, $col).Value = $v
    $col = $col + 1
}

# Row 89: Tunez
$rowValues = @(922, 4, 194, 690, 20, 0, 38)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(89, $col).Value = $v
    $col = $col + 1
}

# Row 97: Costa Rica
$rowValues = @(687, 1, 216, 465, 7, 0, 6)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(97, $col).Value = $v
    $col = $col + 1
}

# Row 126: Jamaica
$rowValues = @(257, 5, 28, 222, 0, 1, 7)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(126, $col).Value = $v
    $col = $col + 1
}

# Row 132: Gabon
$rowValues = @(172, 5, 26, 143, 1, 1, 3)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(132, $col).Value = $v
    $col = $col + 1
}

# Row 133: Martinica
$rowValues = @(170, 6, 77, 79, 6, 0, 14)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(133, $col).Value = $v
    $col = $col + 1
}

# Row 155: Uganda
$rowValues = @(75, 1, 46, 29, 0, 0, 0)
$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(155, $col).Value = $v
    $col = $col + 1
}
